$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the three "Noor.Uddin" related rows: change "ignore" flag (column C) from "no" to "yes"
$ws.Range("C5").Value = "yes"
$ws.Range("C6").Value = "yes"
$ws.Range("C7").Value = "yes"

# Append new rows for a new user "priya.giri" plus two related service rows
$ws.Range("A20").Value = "priya.giri"
$ws.Range("B20").Value = "password3"
$ws.Range("C20").Value = "no"
$ws.Range("D20").Value = "PG"

$ws.Range("A21").Value = "Manufacturer78H14_NU"
$ws.Range("B21").Value = "MHRA12345A"
$ws.Range("C21").Value = "no"
$ws.Range("D21").Value = "PG"

$ws.Range("A22").Value = "AuthorisedRep78H14_NU"
$ws.Range("B22").Value = "MHRA12345A"
$ws.Range("C22").Value = "no"
$ws.Range("D22").Value = "PG"

# Apply the same style used by existing username/password cells to the new rows
$ws.Range("A2:B2").Copy()
$ws.Range("A20:B22").PasteSpecial(-4122)

# Update selection to reflect the last edited cell
$ws.Range("A21").Select()
